# 09/08/2018 RUP - DSDM
# Add new requirement UR-006 ("Formulario para el registro") to the
# Requirements sheet and refresh the affected view state.

$wb = $excel.ActiveWorkbook
$dashboard     = $wb.Worksheets.Item("Dashboard")
$requirements  = $wb.Worksheets.Item("Requirements")

# --- New requirement row (Requirements!A24:L24) -----------------------
$requirements.Range("A24").Value = "UR-006"
$requirements.Range("B24").Value = "Formulario para el registro"
$requirements.Range("C24").Value = "Interfaz"
$requirements.Range("D24").Value = "Se requiere que el formulario para el registro y visualización de la información sea personalizable por roles de usuario."
$requirements.Range("F24").Value = "Manuel Quesada"
$requirements.Range("J24").Value = "Nuevo"
$requirements.Range("K24").Value = "Alta"
$requirements.Range("L24").Value = "Alta"

# The new row wraps onto two lines once populated - match the taller row.
$requirements.Rows.Item(24).RowHeight = 30

# Dashboard's summary description row also grew to a taller, wrapped row.
$dashboard.Rows.Item(8).RowHeight = 45

# --- View state ---------------------------------------------------------
# Dashboard's viewport scrolled down; its selection (H39) is unchanged.
$dashboard.Activate()
$dashWin = $excel.ActiveWindow
$dashWin.ScrollRow = 13
$dashWin.ScrollColumn = 1

# Requirements is (and stays) the active sheet/tab; move its viewport
# down a row and update the active selection/cell.
$requirements.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$requirements.Range("A25").Select()
